$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (column D) and volume-change (column E) values.
# Column D values are plain numeric-looking text (e.g. "603.14") that must stay
# as text, matching the original inlineStr cells, so we force a Text number
# format before assigning and then clear the format delta so no stray style
# index is introduced (original cells carry no explicit style).

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '63.936.14'
$cell.ClearFormats()
$ws.Range("E2").Value = '  -1.06%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.137.29'
$cell.ClearFormats()
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("E4").Value = '  -0.04%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '603.14'
$cell.ClearFormats()
$ws.Range("E5").Value = '  -1.58%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '142.97'
$cell.ClearFormats()
$ws.Range("E6").Value = '  -3.44%  '
$ws.Range("E7").Value = '  -0.02%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '3.143.40'
$cell.ClearFormats()
$ws.Range("E8").Value = '  -0.13%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.523'
$cell.ClearFormats()
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("E10").Value = '  -1.60%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '5.40'
$cell.ClearFormats()
$ws.Range("E11").Value = '  -1.93%  '
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("E13").Value = '  -1.36%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '35.03'
$cell.ClearFormats()
$ws.Range("E14").Value = '  -1.88%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '3.653.97'
$cell.ClearFormats()
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("E16").Value = '  +2.53%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '63.968.31'
$cell.ClearFormats()
$ws.Range("E17").Value = '  -0.92%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '3.137.79'
$cell.ClearFormats()
$ws.Range("E18").Value = '  -0.53%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '6.82'
$cell.ClearFormats()
$ws.Range("E19").Value = '  -1.44%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '486.67'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +1.22%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '14.68'
$cell.ClearFormats()
$ws.Range("E21").Value = '  +0.03%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.708'
$cell.ClearFormats()
$ws.Range("E22").Value = '  -1.22%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '7.64'
$cell.ClearFormats()
$ws.Range("E23").Value = '  -4.46%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '87.07'
$cell.ClearFormats()
$ws.Range("E24").Value = '  +3.67%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '13.42'
$cell.ClearFormats()
$ws.Range("E26").Value = '  +0.01%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '2.75'
$cell.ClearFormats()
$ws.Range("E27").Value = '  -2.69%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '8.24'
$cell.ClearFormats()
$ws.Range("E28").Value = '  -3.42%  '
$ws.Range("E29").Value = '  -0.50%  '
$ws.Range("E30").Value = '  -1.30%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '27.14'
$cell.ClearFormats()
$ws.Range("E31").Value = '  +2.58%  '
$ws.Range("E32").Value = '  -6.88%  '
$ws.Range("E33").Value = '  -0.03%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '2.64'
$cell.ClearFormats()
$ws.Range("E34").Value = '  -3.27%  '
$ws.Range("E35").Value = '  -2.99%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '6.02'
$cell.ClearFormats()
$ws.Range("E36").Value = '  +0.08%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '52.60'
$cell.ClearFormats()
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("E38").Value = '  -5.56%  '
$ws.Range("E39").Value = '  -7.23%  '
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '437.50'
$cell.ClearFormats()
$ws.Range("E40").Value = '  -4.63%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.0396'
$cell.ClearFormats()
$ws.Range("E41").Value = '  -1.28%  '
$ws.Range("E42").Value = '  +0.19%  '
$ws.Range("E43").Value = '  -1.36%  '
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '2.881.74'
$cell.ClearFormats()
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("E45").Value = '  -3.16%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.20'
$cell.ClearFormats()
$ws.Range("E46").Value = '  -5.35%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '2.39'
$cell.ClearFormats()
$ws.Range("E47").Value = '  -3.02%  '
$ws.Range("E48").Value = '  -0.08%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '25.93'
$cell.ClearFormats()
$ws.Range("E49").Value = '  -2.19%  '
$ws.Range("E50").Value = '  -0.27%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '121.02'
$cell.ClearFormats()
$ws.Range("E51").Value = '  +0.42%  '
